$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update labels in column A
$ws.Range("A2").Value = "Anzahl Kinder unter 25"
$ws.Range("A4").Value = "Beitragsbemessungsgrenze PV "
$ws.Range("A5").Value = "Jahresarbeitsentgeltgrenze PV "

# Update entry date value (stored as text)
$ws.Range("B6").Value = "01.01.2025"

# Update selection to match new view state
$ws.Range("B7").Select()
